# Swap the species-observation data between row 2 and row 3 for the
# columns that actually differ between the two records (A, B, E, F, G, H,
# Q, R, Z, AB). The remaining columns already hold identical values in
# both rows, so only these need to be exchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("A", "B", "E", "F", "G", "H", "Q", "R", "Z", "AB")

foreach ($col in $cols) {
    $cell2 = $ws.Range($col + "2")
    $cell3 = $ws.Range($col + "3")

    $val2 = $cell2.Value2
    $val3 = $cell3.Value2

    $cell2.Value2 = $val3
    $cell3.Value2 = $val2
}
